# Report2.xlsx rework:
#  - remove the sort state Excel had stored on the sheet
#  - drop the "Min Mark" / "Max Mark" columns (D:E)
#  - fix the "Session" header typo -> "Seeion", rename "Group" -> "Examainer"
#  - replace the IP-xx group labels with examiner names (2 plain, 2 wrapped w/ leading line break)
#  - reorder/replace the data rows, add a 5th data row
#  - give the two multi-line name cells wrapped text + a taller row

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the sort state stored on the sheet (was sorted by the now-removed "Group" column).
$ws.Sort.SortFields.Clear()

# Remove the Min Mark / Max Mark columns entirely (shifts nothing after them).
$ws.Range("D1:E1").EntireColumn.Delete()

# --- Header row ---
$ws.Range("A1").Value = "Seeion"
$ws.Range("B1").Value = "Examainer"
$ws.Range("C1").Value = "Average Mark"

# --- Data rows ---
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = "Grishchenko Gennady Pavlovich"
$ws.Range("C2").Value = 5.5

$ws.Range("A3").Value = 2
$ws.Range("B3").Value = "Fedorov Nikolay Vladimirovich"
$ws.Range("C3").Value = 7

$ws.Range("A4").Value = 3
$ws.Range("B4").Value = "`n`nBelov Denis Khristoforovich"
$ws.Range("C4").Value = 3.5

$ws.Range("A5").Value = 3
$ws.Range("B5").Value = "`n`nPavlov Ruslan Artemovich"
$ws.Range("C5").Value = 5

# Wrap text + taller rows for the two multi-line examiner names.
$ws.Range("B4").WrapText = $true
$ws.Range("B5").WrapText = $true
$ws.Rows.Item(4).RowHeight = 72
$ws.Rows.Item(5).RowHeight = 72
